$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 12)
$c.Formula = "=""263470844"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(9, 12)
$c.Formula = "=""263478454"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(10, 12)
$c.Formula = "=""263479841"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(11, 12)
$c.Formula = "=""263481152"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(12, 12)
$c.Formula = "=""263482411"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(13, 12)
$c.Formula = "=""263483258"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(14, 12)
$c.Formula = "=""263484672"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(15, 12)
$c.Formula = "=""263485677"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(16, 12)
$c.Formula = "=""263486869"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(17, 12)
$c.Formula = "=""263487716"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(18, 12)
$c.Formula = "=""263488667"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(19, 12)
$c.Formula = "=""263489735"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(20, 12)
$c.Formula = "=""263490665"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(21, 12)
$c.Formula = "=""263491579"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(22, 12)
$c.Formula = "=""263492575"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(23, 12)
$c.Formula = "=""263493487"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(24, 12)
$c.Formula = "=""263495020"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(25, 12)
$c.Formula = "=""263495977"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(26, 12)
$c.Formula = "=""263497601"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(27, 12)
$c.Formula = "=""263498608"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(28, 12)
$c.Formula = "=""263499900"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(29, 12)
$c.Formula = "=""263500770"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(30, 12)
$c.Formula = "=""263502177"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(31, 12)
$c.Formula = "=""263503395"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(32, 12)
$c.Formula = "=""263504687"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(33, 12)
$c.Formula = "=""263505764"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(34, 12)
$c.Formula = "=""263506898"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(35, 12)
$c.Formula = "=""263507966"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(36, 12)
$c.Formula = "=""263509239"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(37, 12)
$c.Formula = "=""263510354"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(38, 12)
$c.Formula = "=""263511684"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(39, 12)
$c.Formula = "=""263513086"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(40, 12)
$c.Formula = "=""263514187"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(41, 12)
$c.Formula = "=""263515276"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(42, 12)
$c.Formula = "=""263516565"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(43, 12)
$c.Formula = "=""263517778"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(44, 12)
$c.Formula = "=""263519166"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(45, 12)
$c.Formula = "=""263520234"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(46, 12)
$c.Formula = "=""263521315"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(47, 12)
$c.Formula = "=""263522919"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(48, 12)
$c.Formula = "=""263524429"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(49, 12)
$c.Formula = "=""263525839"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(50, 12)
$c.Formula = "=""263528473"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(51, 12)
$c.Formula = "=""263529715"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(52, 12)
$c.Formula = "=""263530954"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(53, 12)
$c.Formula = "=""263527101"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(54, 12)
$c.Formula = "=""263532173"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(55, 12)
$c.Formula = "=""263533606"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(56, 12)
$c.Formula = "=""263534859"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(57, 12)
$c.Formula = "=""263536116"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(58, 12)
$c.Formula = "=""263537894"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(59, 12)
$c.Formula = "=""263539512"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(60, 12)
$c.Formula = "=""263541161"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(61, 12)
$c.Formula = "=""263565290"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(62, 12)
$c.Formula = "=""263568735"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(63, 12)
$c.Formula = "=""263571081"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(64, 12)
$c.Formula = "=""263574398"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(65, 12)
$c.Formula = "=""263576798"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(66, 12)
$c.Formula = "=""263580486"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(67, 12)
$c.Formula = "=""263583431"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(68, 12)
$c.Formula = "=""263586110"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(69, 12)
$c.Formula = "=""263543793"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(70, 12)
$c.Formula = "=""263556898"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(71, 12)
$c.Formula = "=""263558397"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(72, 12)
$c.Formula = "=""263546085"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(73, 12)
$c.Formula = "=""263550218"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(74, 12)
$c.Formula = "=""263553959"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(75, 12)
$c.Formula = "=""263559829"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(76, 12)
$c.Formula = "=""263561477"""
$c.Copy()
$c.PasteSpecial(-4163)

$c = $ws.Cells.Item(77, 12)
$c.Formula = "=""263563419"""
$c.Copy()
$c.PasteSpecial(-4163)
